$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Logs sheet: append row 8 with the new test-mail entry
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A8").Value = "Hebben we EcoPro-700 nog op voorraad?"
$logs.Range("B8").Value = "mailmind.test@zohomail.eu"
$logs.Range("C8").Value = "Testmail #6: Hebben we EcoPro-700 nog op voorraad?"
$logs.Range("D8").Value = "Productinformatie"
$logs.Range("E8").Value = "Beste klant," + [char]10 + "Bedankt voor uw vraag. Op dit moment hebben we EcoPro-700 niet op voorraad. We verwachten binnenkort nieuwe voorraad binnen te krijgen. Mocht u nog vragen hebben of een pre-order willen plaatsen, neem dan gerust contact met ons op." + [char]10 + "Met vriendelijke groet," + [char]10 + "[Bedrijfsnaam]"
$logs.Range("F8").Value = "2025-07-31 21:34:42"
$logs.Range("G8").Value = "Ja"
$logs.Range("H8").Value = "Nee"
$logs.Range("I8").Value = "Ja"
$logs.Range("J8").Value = "Nee"

# Extend the conditional-formatting ranges from row 7 to row 8 (one
# ModifyAppliesToRange call per block updates every cfRule sharing that sqref).
$logs.Range("D2:D7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D8"))
$logs.Range("G2:G7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G8"))
$logs.Range("H2:H7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H8"))
$logs.Range("I2:I7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I8"))
$logs.Range("J2:J7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J8"))

# ---------------------------------------------------------------------------
# 2. Dashboard sheet: append row 5 (Productinformatie / 1)
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A5").Value = "Productinformatie"
$dash.Range("B5").Value = 1

# ---------------------------------------------------------------------------
# 3. Chart on Dashboard sheet: widen the category/value series references
#    from row 4 to row 5.
# ---------------------------------------------------------------------------
$chart = $dash.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$5,'Dashboard'!`$B`$2:`$B`$5,1)"
